$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 18782
$ws1.Range("F22").Value = 7811
$ws1.Range("F26").Value = 1232
$ws1.Range("F34").Value = 5366

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 18782
$ws4.Range("F23").Value = 7811
$ws4.Range("F27").Value = 1232
$ws4.Range("F37").Value = 5366
